# "sea sky round 1"
# Fill in the Round_1 (column O, top bracket mirror) winners for the first
# eight match-ups with the winning animal's name - mirroring the values
# already present in the adjacent "key" column P for each row, except for
# row 10 which introduces a brand-new name (Eclectus Parrot) not yet used
# anywhere else in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("O2").Value  = "Orca"
$ws.Range("O4").Value  = "Olive Sea Snake"
$ws.Range("O6").Value  = "Common Map Turtle"
$ws.Range("O8").Value  = "Blanket Octopus"
$ws.Range("O10").Value = "Eclectus Parrot"
$ws.Range("O12").Value = "Steller's Sea Eagle"
$ws.Range("O14").Value = "Indian Fruit Bat"
$ws.Range("O16").Value = "Hawaiian Monk Seal"

# Move the view / selection the way the author left it before saving.
$ws.Activate()
$ws.Range("O5").Select()
